$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 796, pushing the existing row 796 (and everything
# below it) down by one. This mirrors the diff, which inserts a new data
# row ("2026/02/07", "土", 19, 201) right before the former row 796
# ("2026/12/29", ...), shifting the old rows 796-837 down to 797-838.
$ws.Rows(796).Insert()

# Seed the new row by copying the row immediately above it (row 793, which
# already holds the same date "2026/02/07" as inline/shared text rather
# than an auto-converted date serial) so the new cells inherit the correct
# "plain text date" storage and carry no stray number formatting.
$ws.Range("A793:D793").Copy($ws.Range("A796:D796"))

# Now overwrite the time / count / ranking columns with the new row's
# actual values (the date "2026/02/07" and day-of-week "土" already match
# after the copy above).
$ws.Range("B796").Value = "土"
$ws.Range("C796").Value = 19
$ws.Range("D796").Value = 201
